$wb = $excel.ActiveWorkbook

$sheetNames = @("Property1", "Property2")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    # Insert a new row above row 8 (shifts old row 8 -> 9, old row 9 -> 10)
    $ws.Rows.Item(8).Insert()
    # Copy formatting/values of the row that is now row 9 (old row 8) into the new row 8
    $ws.Rows.Item(9).Copy()
    $ws.Rows.Item(8).PasteSpecial()
    # Set the label of the new row 8 to "Force"
    $ws.Cells.Item(8, 1).Value = "Force"
}

$excel.CutCopyMode = $false
